$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.29619540223933
$ws.Range("C2").Value = 10.86081673650695
$ws.Range("D2").Value = 5.839207731680802
$ws.Range("E2").Value = 9.411092624575874
$ws.Range("F2").Value = 34.47678547254818
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("M2").Value = 16.72178415589122
$ws.Range("N2").Value = 18.6560943872973

$ws.Range("B3").Value = 16.6677017318
$ws.Range("C3").Value = 10.25130124603327
$ws.Range("D3").Value = 5.855014677497945
$ws.Range("E3").Value = 9.327985292081809
$ws.Range("F3").Value = 33.95775021068333
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("M3").Value = 16.43211586361197
$ws.Range("N3").Value = 18.70455971377591

$ws.Range("B4").Value = 16.27552129864773
$ws.Range("C4").Value = 9.861107621943241
$ws.Range("D4").Value = 5.865738270268666
$ws.Range("E4").Value = 9.279504279467176
$ws.Range("F4").Value = 33.6467760026635
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("M4").Value = 16.25729199141967
$ws.Range("N4").Value = 18.73625885600948

$ws.Range("B5").Value = 16.1144089506207
$ws.Range("C5").Value = 9.698271865777079
$ws.Range("D5").Value = 5.870362626257373
$ws.Range("E5").Value = 9.260404042046821
$ws.Range("F5").Value = 33.52214516453453
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("M5").Value = 16.18691059239602
$ws.Range("N5").Value = 18.7496635547037

$ws.Range("B6").Value = 16.08758689441258
$ws.Range("C6").Value = 9.671007875203907
$ws.Range("D6").Value = 5.871145816017926
$ws.Range("E6").Value = 9.257272521522617
$ws.Range("F6").Value = 33.50158115043912
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("M6").Value = 16.17527863821542
$ws.Range("N6").Value = 18.75191877925722

$ws.Range("B7").Value = 16.27335333375519
$ws.Range("C7").Value = 9.858926798303209
$ws.Range("D7").Value = 5.865799607966428
$ws.Range("E7").Value = 9.279244010580337
$ws.Range("F7").Value = 33.64508651110076
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("M7").Value = 16.25633918962999
$ws.Range("N7").Value = 18.73643766586052

$ws.Range("B8").Value = 17.08095773709594
$ws.Range("C8").Value = 10.65406455213084
$ws.Range("D8").Value = 5.844445704137271
$ws.Range("E8").Value = 9.381918366211991
$ws.Range("F8").Value = 34.29631427241144
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("M8").Value = 16.62133595036878
$ws.Range("N8").Value = 18.67240166067977

$ws.Range("B9").Value = 18.60358081520028
$ws.Range("C9").Value = 12.08040843514816
$ws.Range("D9").Value = 5.810720801783704
$ws.Range("E9").Value = 9.602719317891172
$ws.Range("F9").Value = 35.62764776307567
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("M9").Value = 17.35683165452803
$ws.Range("N9").Value = 18.56227980556225

$ws.Range("B10").Value = 19.67171448310222
$ws.Range("C10").Value = 13.04064062650107
$ws.Range("D10").Value = 5.791011652195832
$ws.Range("E10").Value = 9.775709390113844
$ws.Range("F10").Value = 36.62942440421095
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("M10").Value = 17.90346753237175
$ws.Range("N10").Value = 18.49086276867231

$ws.Range("B11").Value = 20.14433384215484
$ws.Range("C11").Value = 13.45745444288346
$ws.Range("D11").Value = 5.783167833827736
$ws.Range("E11").Value = 9.85649159793625
$ws.Range("F11").Value = 37.08828772162438
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("M11").Value = 18.15231892267457
$ws.Range("N11").Value = 18.46044964086826

$ws.Range("B12").Value = 20.32122904065497
$ws.Range("C12").Value = 13.61235513597234
$ws.Range("D12").Value = 5.780360740087978
$ws.Range("E12").Value = 9.887359579655083
$ws.Range("F12").Value = 37.26233564311804
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("M12").Value = 18.2464837215477
$ws.Range("N12").Value = 18.44923278991249

$ws.Range("B13").Value = 20.28322622476476
$ws.Range("C13").Value = 13.57912599240407
$ws.Range("D13").Value = 5.780958010379435
$ws.Range("E13").Value = 9.880699632615233
$ws.Range("F13").Value = 37.22484118416079
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("M13").Value = 18.22620824482179
$ws.Range("N13").Value = 18.45163517374428

$ws.Range("B14").Value = 20.15892958924004
$ws.Range("C14").Value = 13.47025739628444
$ws.Range("D14").Value = 5.782933611309648
$ws.Range("E14").Value = 9.85902569571436
$ws.Range("F14").Value = 37.10260195040468
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("M14").Value = 18.16006785458023
$ws.Range("N14").Value = 18.45952080329776

$ws.Range("B15").Value = 20.08251945131364
$ws.Range("C15").Value = 13.40318807055225
$ws.Range("D15").Value = 5.784165031281183
$ws.Range("E15").Value = 9.845785239874767
$ws.Range("F15").Value = 37.02775930226785
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("M15").Value = 18.11954307031186
$ws.Range("N15").Value = 18.46439008800886

$ws.Range("B16").Value = 19.64054581475941
$ws.Range("C16").Value = 13.01299282073061
$ws.Range("D16").Value = 5.79154699909013
$ws.Range("E16").Value = 9.770470168145781
$ws.Range("F16").Value = 36.59948566380621
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("M16").Value = 17.8871998958947
$ws.Range("N16").Value = 18.49289222311604

$ws.Range("B17").Value = 19.36588329044428
$ws.Range("C17").Value = 12.76844976273585
$ws.Range("D17").Value = 5.796364274225977
$ws.Range("E17").Value = 9.724785616391507
$ws.Range("F17").Value = 36.33744074762071
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("M17").Value = 17.74464234408298
$ws.Range("N17").Value = 18.510909884177

$ws.Range("B18").Value = 19.20666322437895
$ws.Range("C18").Value = 12.62591522515742
$ws.Range("D18").Value = 5.799240537232922
$ws.Range("E18").Value = 9.698707019681956
$ws.Range("F18").Value = 36.1870257163837
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("M18").Value = 17.66266987233013
$ws.Range("N18").Value = 18.52146836234498

$ws.Range("B19").Value = 19.152546223201
$ws.Range("C19").Value = 12.57733471984115
$ws.Range("D19").Value = 5.800232452911382
$ws.Range("E19").Value = 9.68991193764343
$ws.Range("F19").Value = 36.13615551252162
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("M19").Value = 17.63492232240064
$ws.Range("N19").Value = 18.52507675937433

$ws.Range("B20").Value = 19.39525133752902
$ws.Range("C20").Value = 12.79467683423342
$ws.Range("D20").Value = 5.795840535497616
$ws.Range("E20").Value = 9.729628491173735
$ws.Range("F20").Value = 36.36530539749575
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("M20").Value = 17.75981618209854
$ws.Range("N20").Value = 18.50897166007399

$ws.Range("B21").Value = 20.19549607187576
$ws.Range("C21").Value = 13.50231489355206
$ws.Range("D21").Value = 5.782348885664971
$ws.Range("E21").Value = 9.865384508578869
$ws.Range("F21").Value = 37.1385001133826
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("M21").Value = 18.17949753376438
$ws.Range("N21").Value = 18.45719644957568

$ws.Range("B22").Value = 20.70634198437685
$ws.Range("C22").Value = 13.9476547921175
$ws.Range("D22").Value = 5.774483378245203
$ws.Range("E22").Value = 9.955715605511656
$ws.Range("F22").Value = 37.6454223416691
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("M22").Value = 18.45333472478572
$ws.Range("N22").Value = 18.42510741093441

$ws.Range("B23").Value = 20.43485610602161
$ws.Range("C23").Value = 13.71155403406794
$ws.Range("D23").Value = 5.778593617984659
$ws.Range("E23").Value = 9.907364824328354
$ws.Range("F23").Value = 37.37477644998115
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("M23").Value = 18.30725452786946
$ws.Range("N23").Value = 18.44207335981083

$ws.Range("B24").Value = 19.38197812009781
$ws.Range("C24").Value = 12.78282562389483
$ws.Range("D24").Value = 5.796076985422765
$ws.Range("E24").Value = 9.727438446226767
$ws.Range("F24").Value = 36.35270703503137
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("M24").Value = 17.7529561269282
$ws.Range("N24").Value = 18.50984730871384

$ws.Range("B25").Value = 18.19967843479009
$ws.Range("C25").Value = 11.70955322661465
$ws.Range("D25").Value = 5.818961487626365
$ws.Range("E25").Value = 9.541012466094244
$ws.Range("F25").Value = 35.26264829728986
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("M25").Value = 17.15634805218522
$ws.Range("N25").Value = 18.59040941440993
